$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.550.13"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.755.18"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'324.45"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'0.4571"
$ws.Range("E7").Value = "  +2.33%  "
$ws.Range("D8").Value = "'0.3559"
$ws.Range("E8").Value = "  -1.84%  "
$ws.Range("D9").Value = "'0.07448"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").Value = "'41.50"
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("D11").Value = "'1.085"
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "'20.77"
$ws.Range("D14").Value = "'6.001"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").Value = "1.748.45"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "'93.57"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").Value = "'0.06402"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").Value = "'5.736"
$ws.Range("E22").Value = "  -2.13%  "
$ws.Range("D23").Value = "27.592.08"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "'2.070"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("D26").Value = "'164.87"
$ws.Range("E26").Value = "  +1.94%  "
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").Value = "1.952.84"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "'2.127"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "'125.20"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "'1.078"
$ws.Range("D32").Value = "'0.09229"
$ws.Range("E32").Value = "  +2.37%  "
$ws.Range("D33").Value = "'3.661"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("D36").Value = "'0.02277"
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("D37").Value = "'0.2087"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("E39").Value = "  -1.41%  "
$ws.Range("D40").Value = "'4.920"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").Value = "'1.181"
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").Value = "'1.384"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").Value = "'7.760"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "'13.18"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").Value = "'3.714"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "'0.5865"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").Value = "'121.89"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("D48").Value = "'1.933"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("D49").Value = "'0.06894"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("D51").Value = "'71.97"
$ws.Range("E51").Value = "  -0.64%  "
